$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = "João Rodrigues-Desenho Técnico"
$ws.Range("D3").Value = "Euclides-Gestão"
$ws.Range("E3").Value = "José Ferreira-Tecnologia dos Materiais"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "Euclides-Gestão"
$ws.Range("E4").Value = "José Ferreira-Tecnologia dos Materiais"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("C6").Value = "Andre Lucca-Circuitos Elétricos"
$ws.Range("D6").Value = "-"
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("B7").Value = "João Rodrigues-Desenho Técnico"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "Andre Lucca-Circuitos Elétricos"
